$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the G column (Week topics) text for rows 4-15 -----------------
# These assignments cause the old, now-unreferenced shared strings
# ("Introduction to  React JS and Material UI." and
#  "MuiTypography,MuiTextbox,MuiButton") to be dropped, and the new /
# reshuffled strings to be appended in the same order the diff expects.
$ws.Range("G4").Value2  = "Introduction to  React JS,Functional Component and Arrow Function"
$ws.Range("G5").Value2  = "Introduction to  Material UI,MuiTypography,MuiTextbox,MuiButton"
$ws.Range("G6").Value2  = "MuiSwitch,MuiRadibutton,MuiCheckbox"
$ws.Range("G7").Value2  = "MuiSelect,MuiAutocomplete"
$ws.Range("G8").Value2  = "MuiRating,MuiCard,MuiLayout"
$ws.Range("G9").Value2  = "MuiAccordion,MuiImageList,MuiNavbar"
$ws.Range("G10").Value2 = "MuiLink,MuiBreadcrumbs,MuiDrawer"
$ws.Range("G11").Value2 = "MuiSpeedDial,MuiAvatar,MuiBadge"
$ws.Range("G12").Value2 = "MuiList,MuiResponsiveness,MuiSkeleton"
$ws.Range("G13").Value2 = "MuiCarousel,Swiperdemo"
$ws.Range("G14").Value2 = "React Router"
$ws.Range("G15").Value2 = "Sample Programs to revise of above topics"

# --- Widen column G so the longer text fits (bestFit-style width) ---------
$ws.Columns("G").ColumnWidth = 64.8

# --- Remove the now-blank spacer row (old row 18); rows below shift up ----
$ws.Rows("18:18").Delete()

# --- Restore the view: select the new active cell -------------------------
[void]$ws.Range("G19").Select()

Write-Output "Edit complete"
